$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Valor Mora" header date value (row 11)
$ws.Range("E11").Value = 35711

# Delete entire row 17 (second Daniel Brieva Meza entry); rows below shift up
$ws.Rows("17").Delete()

# Row 16 now becomes the updated entry for Andres Felipe Ponce Morales
$ws.Range("C16").Value = "1050970745"
$ws.Range("D16").Value = "ANDRES FELIPE PONCE MORALES"
$ws.Range("E16").Value = "2309"
$ws.Range("F16").Value = 1547
$ws.Range("G16").Value = 1423500

# Row 17 (previously row 18) becomes the new entry for Darwin Reales Castro
$ws.Range("C17").Value = "1002059825"
$ws.Range("D17").Value = "DARWIN REALES CASTRO"
$ws.Range("E17").Value = "2508"
$ws.Range("F17").Value = 34164
$ws.Range("G17").Value = 1423500
